$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.388106333333333
$ws.Range("H2").Value = 28.164319
$ws.Range("I2").Value = 0.2414596449149976
$ws.Range("J2").Value = 0.2414596449149975
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 104.6477755292914
$ws.Range("R2").Value = 941.8299797636221
$ws.Range("S2").Value = 0.06265409560696686
$ws.Range("T2").Value = 0.06265409560696684

# Row 3
$ws.Range("G3").Value = 9.388106333333333
$ws.Range("H3").Value = 28.164319
$ws.Range("I3").Value = 0.2414596449149976
$ws.Range("J3").Value = 0.2414596449149975
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 259.9143897271809
$ws.Range("R3").Value = 2339.229507544629
$ws.Range("S3").Value = 0.1556144021335177
$ws.Range("T3").Value = 0.1556144021335177

# Row 4
$ws.Range("G4").Value = 9.388106333333333
$ws.Range("H4").Value = 28.164319
$ws.Range("I4").Value = 0.2414596449149976
$ws.Range("J4").Value = 0.2414596449149975
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 38.73492930149856
$ws.Range("R4").Value = 348.614363713487
$ws.Range("S4").Value = 0.02319114717451297
$ws.Range("T4").Value = 0.02319114717451297

# Row 5
$ws.Range("I5").Value = 0.5770971896641285
$ws.Range("J5").Value = 0.5770971896641284
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 250.11192733186
$ws.Range("R5").Value = 2251.00734598674
$ws.Range("S5").Value = 0.1497455299764767
$ws.Range("T5").Value = 0.1497455299764767

# Row 6
$ws.Range("I6").Value = 0.5770971896641285
$ws.Range("J6").Value = 0.5770971896641284
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.3719239882678164
$ws.Range("T6").Value = 0.3719239882678163

# Row 7
$ws.Range("I7").Value = 0.5770971896641285
$ws.Range("J7").Value = 0.5770971896641284
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.05542767141983546
$ws.Range("T7").Value = 0.05542767141983545

# Row 8
$ws.Range("I8").Value = 0.181443165420874
$ws.Range("J8").Value = 0.1814431654208739
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 78.63684075644201
$ws.Range("R8").Value = 707.7315668079781
$ws.Range("S8").Value = 0.04708098298376998
$ws.Range("T8").Value = 0.04708098298376997

# Row 9
$ws.Range("I9").Value = 0.181443165420874
$ws.Range("J9").Value = 0.1814431654208739
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.1169353567057636
$ws.Range("T9").Value = 0.1169353567057636

# Row 10
$ws.Range("I10").Value = 0.181443165420874
$ws.Range("J10").Value = 0.1814431654208739
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.01742682573134038
$ws.Range("T10").Value = 0.01742682573134038
